$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 100
$ws_ALC.Range("H100").Value = 911.1667
$ws_ALC.Range("I100").Value = 799.625
$ws_ALC.Range("J100").Value = 1134.25
$ws_ALC.Range("K100").Value = 799.625
$ws_ALC.Range("L100").Value = 1134.25
$ws_ALC.Range("M100").Value = -258.625
$ws_ALC.Range("N100").Value = -2216.25

# ALC row 106
$ws_ALC.Range("H106").Value = 1567.5
$ws_ALC.Range("I106").Value = 1376.25
$ws_ALC.Range("J106").Value = 1950
$ws_ALC.Range("K106").Value = 1376.25
$ws_ALC.Range("L106").Value = 1950
$ws_ALC.Range("M106").Value = -745.25
$ws_ALC.Range("N106").Value = -3212

# ALC row 137
$ws_ALC.Range("H137").Value = 2580.4285
$ws_ALC.Range("I137").Value = 2095.9092
$ws_ALC.Range("J137").Value = 3400.3845
$ws_ALC.Range("K137").Value = 6287.7276
$ws_ALC.Range("L137").Value = 10201.1535
$ws_ALC.Range("M137").Value = -3737.7276
$ws_ALC.Range("N137").Value = -15301.1535

# ALC row 138
$ws_ALC.Range("H138").Value = 3367.2727
$ws_ALC.Range("I138").Value = 1535.7812
$ws_ALC.Range("J138").Value = 5091.0293
$ws_ALC.Range("K138").Value = 4607.3436
$ws_ALC.Range("L138").Value = 15273.0879
$ws_ALC.Range("M138").Value = 532.6563999999998
$ws_ALC.Range("N138").Value = -25553.0879

# ARM row 2
$ws_ARM.Range("H2").Value = 1292.7084
$ws_ARM.Range("I2").Value = 1536.3077
$ws_ARM.Range("J2").Value = 1004.8182
$ws_ARM.Range("K2").Value = 1536.3077
$ws_ARM.Range("L2").Value = 1004.8182
$ws_ARM.Range("M2").Value = -1423.3077
$ws_ARM.Range("N2").Value = -1230.8182

# ARM row 32
$ws_ARM.Range("H32").Value = 17303.867
$ws_ARM.Range("I32").Value = 17893.434
$ws_ARM.Range("J32").Value = 12366.25
$ws_ARM.Range("K32").Value = 17893.434
$ws_ARM.Range("L32").Value = 12366.25
$ws_ARM.Range("M32").Value = -17606.434
$ws_ARM.Range("N32").Value = -12940.25

# ARM row 34
$ws_ARM.Range("H34").Value = 20028
$ws_ARM.Range("J34").Value = 20028
$ws_ARM.Range("L34").Value = 20028
$ws_ARM.Range("N34").Value = -20570

# ARM row 116
$ws_ARM.Range("H116").Value = 1292.7084
$ws_ARM.Range("I116").Value = 1536.3077
$ws_ARM.Range("J116").Value = 1004.8182
$ws_ARM.Range("K116").Value = 1536.3077
$ws_ARM.Range("L116").Value = 1004.8182
$ws_ARM.Range("M116").Value = 757.6922999999999
$ws_ARM.Range("N116").Value = -5592.8182

# ARM row 132
$ws_ARM.Range("H132").Value = 6409.3105
$ws_ARM.Range("I132").Value = 2398.923
$ws_ARM.Range("J132").Value = 9667.75
$ws_ARM.Range("K132").Value = 7196.768999999999
$ws_ARM.Range("L132").Value = 29003.25
$ws_ARM.Range("M132").Value = -4666.768999999999
$ws_ARM.Range("N132").Value = -34063.25

# BSM row 3
$ws_BSM.Range("H3").Value = 1292.7084
$ws_BSM.Range("I3").Value = 1536.3077
$ws_BSM.Range("J3").Value = 1004.8182
$ws_BSM.Range("K3").Value = 1536.3077
$ws_BSM.Range("L3").Value = 1004.8182
$ws_BSM.Range("M3").Value = -1422.3077
$ws_BSM.Range("N3").Value = -1232.8182

# BSM row 117
$ws_BSM.Range("H117").Value = 0
$ws_BSM.Range("J117").Value = 0
$ws_BSM.Range("L117").Value = 0
$ws_BSM.Range("N117").ClearContents() | Out-Null

# BSM row 134
$ws_BSM.Range("H134").Value = 19161.19
$ws_BSM.Range("I134").Value = 1845.4878
$ws_BSM.Range("J134").Value = 60922.59
$ws_BSM.Range("K134").Value = 5536.463400000001
$ws_BSM.Range("L134").Value = 182767.77
$ws_BSM.Range("M134").Value = -3001.463400000001
$ws_BSM.Range("N134").Value = -187837.77

# CRP row 16
$ws_CRP.Range("H16").Value = 1031.5769
$ws_CRP.Range("I16").Value = 443.7857
$ws_CRP.Range("K16").Value = 443.7857
$ws_CRP.Range("M16").Value = -156.7857

# CRP row 113
$ws_CRP.Range("H113").Value = 1031.5769
$ws_CRP.Range("I113").Value = 443.7857
$ws_CRP.Range("K113").Value = 443.7857
$ws_CRP.Range("M113").Value = 1726.2143

# CUL row 8
$ws_CUL.Range("H8").Value = 63.46154
$ws_CUL.Range("I8").Value = 63.46154
$ws_CUL.Range("K8").Value = 190.38462
$ws_CUL.Range("M8").Value = -51.38461999999998

# CUL row 92
$ws_CUL.Range("H92").Value = 704.4545000000001
$ws_CUL.Range("J92").Value = 729.5714
$ws_CUL.Range("L92").Value = 2188.7142
$ws_CUL.Range("N92").Value = -4684.7142

# CUL row 131
$ws_CUL.Range("H131").Value = 22946.818
$ws_CUL.Range("I131").Value = 2758.3333
$ws_CUL.Range("J131").Value = 26134.475
$ws_CUL.Range("K131").Value = 8274.999899999999
$ws_CUL.Range("L131").Value = 78403.42499999999
$ws_CUL.Range("M131").Value = -3234.999899999999
$ws_CUL.Range("N131").Value = -88483.42499999999

# GSM row 113
$ws_GSM.Range("H113").Value = 1813.7142
$ws_GSM.Range("I113").Value = 1661.3529
$ws_GSM.Range("J113").Value = 1957.6111
$ws_GSM.Range("K113").Value = 1661.3529
$ws_GSM.Range("L113").Value = 1957.6111
$ws_GSM.Range("M113").Value = 508.6470999999999
$ws_GSM.Range("N113").Value = -6297.6111

# GSM row 122
$ws_GSM.Range("H122").Value = 3642.7693
$ws_GSM.Range("I122").Value = 3525.818
$ws_GSM.Range("K122").Value = 10577.454
$ws_GSM.Range("M122").Value = -8127.454000000002

# GSM row 132
$ws_GSM.Range("H132").Value = 6302.148
$ws_GSM.Range("I132").Value = 2875.182
$ws_GSM.Range("J132").Value = 21380.8
$ws_GSM.Range("K132").Value = 8625.545999999998
$ws_GSM.Range("L132").Value = 64142.39999999999
$ws_GSM.Range("M132").Value = -6095.545999999998
$ws_GSM.Range("N132").Value = -69202.39999999999

# LTW row 22
$ws_LTW.Range("H22").Value = 1266.6666
$ws_LTW.Range("I22").Value = 1150
$ws_LTW.Range("J22").Value = 1500
$ws_LTW.Range("K22").Value = 1150
$ws_LTW.Range("L22").Value = 1500
$ws_LTW.Range("M22").Value = -855
$ws_LTW.Range("N22").Value = -2090

# LTW row 27
$ws_LTW.Range("H27").Value = 1266.6666
$ws_LTW.Range("I27").Value = 1150
$ws_LTW.Range("J27").Value = 1500
$ws_LTW.Range("K27").Value = 1150
$ws_LTW.Range("L27").Value = 1500
$ws_LTW.Range("M27").Value = -1043
$ws_LTW.Range("N27").Value = -1714

# LTW row 40
$ws_LTW.Range("H40").Value = 4172.65
$ws_LTW.Range("I40").Value = 3850.1765
$ws_LTW.Range("J40").Value = 6000
$ws_LTW.Range("K40").Value = 3850.1765
$ws_LTW.Range("L40").Value = 6000
$ws_LTW.Range("M40").Value = -3714.1765
$ws_LTW.Range("N40").Value = -6272

# LTW row 61
$ws_LTW.Range("H61").Value = 780362.4399999999
$ws_LTW.Range("I61").Value = 12913.2
$ws_LTW.Range("K61").Value = 12913.2
$ws_LTW.Range("M61").Value = -12711.2

# LTW row 113
$ws_LTW.Range("H113").Value = 780362.4399999999
$ws_LTW.Range("I113").Value = 12913.2
$ws_LTW.Range("K113").Value = 12913.2
$ws_LTW.Range("M113").Value = -10743.2

# LTW row 136
$ws_LTW.Range("H136").Value = 3583.103
$ws_LTW.Range("I136").Value = 2048.6938
$ws_LTW.Range("K136").Value = 6146.0814
$ws_LTW.Range("M136").Value = -3596.0814

# WVR row 113
$ws_WVR.Range("H113").Value = 6892.2354
$ws_WVR.Range("I113").Value = 13065.25
$ws_WVR.Range("J113").Value = 1405.1111
$ws_WVR.Range("K113").Value = 39195.75
$ws_WVR.Range("L113").Value = 4215.3333
$ws_WVR.Range("M113").Value = -37025.75
$ws_WVR.Range("N113").Value = -8555.3333

# WVR row 122
$ws_WVR.Range("H122").Value = 2534.5676
$ws_WVR.Range("I122").Value = 2015.129
$ws_WVR.Range("J122").Value = 5218.3335
$ws_WVR.Range("K122").Value = 6045.387
$ws_WVR.Range("L122").Value = 15655.0005
$ws_WVR.Range("M122").Value = -3595.387
$ws_WVR.Range("N122").Value = -20555.0005
